# Generate Report for Handoff
# A new handoff just occurred for the "268b50a0-f412-4f69-99e2-079bdfdf1585" file
# (row 4 on both the "zh-cn" and "de-de" per-locale sheets). Update the
# "Latest Handoff Datetime" (column D) on that row to reflect the new handoff
# timestamp for each locale.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("D4").Value = "2016-01-25 06:10:10"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("D4").Value = "2016-01-25 06:10:21"
